$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2367.9048
$ws.Range("J40").Value = 2618.25
$ws.Range("L40").Value = 2618.25
$ws.Range("N40").Value = -2968.25

# Row 64
$ws.Range("H64").Value = 5500
$ws.Range("J64").Value = 5500
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996

# Row 67
$ws.Range("H67").Value = 5500
$ws.Range("J67").Value = 5500
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 2477.842
$ws.Range("I74").Value = 1735.9286
$ws.Range("J74").Value = 4555.2
$ws.Range("K74").Value = 1735.9286
$ws.Range("L74").Value = 4555.2
$ws.Range("M74").Value = -861.9286
$ws.Range("N74").Value = -6303.2

# Row 77
$ws.Range("H77").Value = 2477.842
$ws.Range("I77").Value = 1735.9286
$ws.Range("J77").Value = 4555.2
$ws.Range("K77").Value = 8679.643
$ws.Range("L77").Value = 22776
$ws.Range("M77").Value = -4311.643
$ws.Range("N77").Value = -31512

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 831.375
$ws.Range("I22").Value = 758.1667
$ws.Range("K22").Value = 758.1667
$ws.Range("M22").Value = -585.1667

# Row 64
$ws.Range("H64").Value = 487.83334
$ws.Range("J64").Value = 457.5
$ws.Range("L64").Value = 457.5
$ws.Range("N64").Value = -907.5

# Row 67
$ws.Range("H67").Value = 487.83334
$ws.Range("J67").Value = 457.5
$ws.Range("L67").Value = 457.5
$ws.Range("N67").Value = -2017.5

# Row 105
$ws.Range("H105").Value = 3029.2144
$ws.Range("I105").Value = 2793.7273
$ws.Range("J105").Value = 3892.6667
$ws.Range("K105").Value = 2793.7273
$ws.Range("L105").Value = 3892.6667
$ws.Range("M105").Value = -1046.7273
$ws.Range("N105").Value = -7386.6667

# Row 122
$ws.Range("H122").Value = 85389.5
$ws.Range("J122").Value = 85389.5
$ws.Range("L122").Value = 85389.5
$ws.Range("N122").Value = -95189.5

# Row 134
$ws.Range("H134").Value = 1619.8334
$ws.Range("I134").Value = 430.25
$ws.Range("K134").Value = 1290.75
$ws.Range("M134").Value = 1244.25

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1572.5555
$ws.Range("I16").Value = 1267.5625
$ws.Range("K16").Value = 1267.5625
$ws.Range("M16").Value = -980.5625

# Row 33
$ws.Range("H33").Value = 1499.75
$ws.Range("I33").Value = 1499.75
$ws.Range("K33").Value = 1499.75
$ws.Range("M33").Value = -1120.75

# Row 92
$ws.Range("H92").Value = 50198.668
$ws.Range("J92").Value = 50198.668
$ws.Range("L92").Value = 50198.668
$ws.Range("N92").Value = -55190.668

# Row 105
$ws.Range("H105").Value = 3828.4285
$ws.Range("I105").Value = 3272.6365
$ws.Range("J105").Value = 4439.8
$ws.Range("K105").Value = 3272.6365
$ws.Range("L105").Value = 4439.8
$ws.Range("M105").Value = -1525.6365
$ws.Range("N105").Value = -7933.8

# Row 107
$ws.Range("H107").Value = 590.7
$ws.Range("I107").Value = 403.7143
$ws.Range("K107").Value = 403.7143
$ws.Range("M107").Value = 1516.2857

# Row 113
$ws.Range("H113").Value = 1572.5555
$ws.Range("I113").Value = 1267.5625
$ws.Range("K113").Value = 1267.5625
$ws.Range("M113").Value = 902.4375

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 9999
$ws.Range("I3").Value = 9999
$ws.Range("K3").Value = 29997
$ws.Range("M3").Value = -29885

# Row 18
$ws.Range("H18").Value = 567
$ws.Range("I18").Value = 612.2857
$ws.Range("J18").Value = 250
$ws.Range("K18").Value = 1836.8571
$ws.Range("L18").Value = 750
$ws.Range("M18").Value = -1667.8571
$ws.Range("N18").Value = -1088

# Row 68
$ws.Range("H68").Value = 12506300
$ws.Range("I68").Value = 2749.5
$ws.Range("K68").Value = 8248.5
$ws.Range("M68").Value = -7437.5

# Row 71
$ws.Range("H71").Value = 12506300
$ws.Range("I71").Value = 2749.5
$ws.Range("K71").Value = 24745.5
$ws.Range("M71").Value = -20689.5

# Row 75
$ws.Range("H75").Value = 3668.25
$ws.Range("J75").Value = 3847.7273
$ws.Range("L75").Value = 11543.1819
$ws.Range("N75").Value = -13539.1819

# Row 78
$ws.Range("H78").Value = 3668.25
$ws.Range("J78").Value = 3847.7273
$ws.Range("L78").Value = 34629.5457
$ws.Range("N78").Value = -44613.5457

# Row 130
$ws.Range("H130").Value = 2746.75
$ws.Range("J130").Value = 3192.8
$ws.Range("L130").Value = 9578.400000000001
$ws.Range("N130").Value = -19618.4

# Row 132
$ws.Range("H132").Value = 3766.8845
$ws.Range("J132").Value = 5126.3335
$ws.Range("L132").Value = 46137.0015
$ws.Range("N132").Value = -51197.0015

# Row 134
$ws.Range("H134").Value = 15635.833
$ws.Range("J134").Value = 18363
$ws.Range("L134").Value = 55089
$ws.Range("N134").Value = -65229

# Row 138
$ws.Range("H138").Value = 3060
$ws.Range("I138").Value = 3060
$ws.Range("K138").Value = 9180
$ws.Range("M138").Value = -4040

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2227.3845
$ws.Range("I122").Value = 1371.5
$ws.Range("J122").Value = 5080.3335
$ws.Range("K122").Value = 4114.5
$ws.Range("L122").Value = 15241.0005
$ws.Range("M122").Value = -1664.5
$ws.Range("N122").Value = -20141.0005

# Row 126
$ws.Range("H126").Value = 3010.1667
$ws.Range("I126").Value = 2632
$ws.Range("K126").Value = 7896
$ws.Range("M126").Value = -5426

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8194.549999999999
$ws.Range("I7").Value = 4302.5
$ws.Range("K7").Value = 4302.5
$ws.Range("M7").Value = -4190.5

# Row 22
$ws.Range("H22").Value = 1008.6429
$ws.Range("I22").Value = 864.6667
$ws.Range("J22").Value = 1267.8
$ws.Range("K22").Value = 864.6667
$ws.Range("L22").Value = 1267.8
$ws.Range("M22").Value = -569.6667
$ws.Range("N22").Value = -1857.8

# Row 27
$ws.Range("H27").Value = 1008.6429
$ws.Range("I27").Value = 864.6667
$ws.Range("J27").Value = 1267.8
$ws.Range("K27").Value = 864.6667
$ws.Range("L27").Value = 1267.8
$ws.Range("M27").Value = -757.6667
$ws.Range("N27").Value = -1481.8

# Row 126
$ws.Range("H126").Value = 8194.549999999999
$ws.Range("I126").Value = 4302.5
$ws.Range("K126").Value = 12907.5
$ws.Range("M126").Value = -10437.5

$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 10550
$ws.Range("I74").Value = 10550
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 10550
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -9614
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 10550
$ws.Range("I77").Value = 10550
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 31650
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -26970
$ws.Range("N77").ClearContents()
